# "Processed Salesforce Current Translations":
#   - insert a new "styles" column (B) ahead of the locale columns
#   - insert a new "Attn" row (2) ahead of the existing label rows
# Both inserts shift the previously-existing data right/down, matching the
# unified diff (old column B.. becomes C.., old row 2.. becomes row 3..).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column B ("styles"), existing zh_CN..sheet columns shift to C..T.
$ws.Columns.Item(2).Insert()

# New row 2 ("Attn"), existing Delivery..Ship To rows shift to 3..8.
$ws.Rows.Item(2).Insert()

# Header label for the newly inserted column.
$ws.Range("B1").Value = 'styles'

# Populate the newly inserted "Attn" row across all 20 columns.
$attn = @(
    'Attn',
    'FFF2CC',
    '注意',
    'À l''attention de',
    'Beachtung',
    'Προσοχή',
    'Attenzione',
    '주목',
    'Uwaga',
    'Atenção',
    'Atenção',
    'Atención',
    'Atención',
    'ความสนใจ',
    'Liên Hệ',
    'CustomLabel$CEC_Attention',
    'CustomLabel',
    'Attn',
    'IVP Contact Creation Translations.xlsx',
    'Sheet1'
)
$arr = New-Object 'object[,]' 1,$attn.Length
for ($i = 0; $i -lt $attn.Length; $i++) {
    $arr[0,$i] = $attn[$i]
}
$ws.Range("A2:T2").Value = $arr
